$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for columns D and E so numeric-looking strings
# (e.g. "202.15", "0.999") are not silently converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '67.878.20'
$ws.Range("E2").Value = '  +2.10%  '
$ws.Range("D3").Value = '3.589.66'
$ws.Range("E3").Value = '  +0.66%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").Value = '202.15'
$ws.Range("E5").Value = '  +8.71%  '
$ws.Range("D6").Value = '564.13'
$ws.Range("E6").Value = '  -3.76%  '
$ws.Range("D7").Value = '3.583.04'
$ws.Range("E7").Value = '  +0.65%  '
$ws.Range("D8").Value = '0.619'
$ws.Range("E8").Value = '  +0.88%  '
$ws.Range("E9").Value = '  -0.06%  '
$ws.Range("D10").Value = '0.668'
$ws.Range("E10").Value = '  -0.46%  '
$ws.Range("D11").Value = '60.49'
$ws.Range("E11").Value = '  +13.35%  '
$ws.Range("D12").Value = '0.150'
$ws.Range("E12").Value = '  +2.89%  '
$ws.Range("D13").Value = '0.0000283'
$ws.Range("E13").Value = '  +9.00%  '
$ws.Range("D14").Value = '9.95'
$ws.Range("E14").Value = '  +1.15%  '
$ws.Range("D15").Value = '4.178.28'
$ws.Range("E15").Value = '  +0.90%  '
$ws.Range("D16").Value = '3.598.73'
$ws.Range("E16").Value = '  +0.87%  '
$ws.Range("E17").Value = '  +0.52%  '
$ws.Range("D18").Value = '18.87'
$ws.Range("E18").Value = '  +2.85%  '
$ws.Range("D19").Value = '67.764.48'
$ws.Range("E19").Value = '  +2.06%  '
$ws.Range("D20").Value = '12.32'
$ws.Range("E20").Value = '  +0.65%  '
$ws.Range("D21").Value = '1.07'
$ws.Range("E21").Value = '  +1.55%  '
$ws.Range("D22").Value = '400.40'
$ws.Range("E22").Value = '  +1.15%  '
$ws.Range("D23").Value = '12.96'
$ws.Range("E23").Value = '  +14.54%  '
$ws.Range("D24").Value = '4.15'
$ws.Range("E24").Value = '  -4.92%  '
$ws.Range("D25").Value = '85.04'
$ws.Range("E25").Value = '  -1.27%  '
$ws.Range("D26").Value = '2.90'
$ws.Range("E26").Value = '  -0.38%  '
$ws.Range("D27").Value = '12.56'
$ws.Range("E27").Value = '  +0.59%  '
$ws.Range("D28").Value = '3.93'
$ws.Range("E28").Value = '  +10.72%  '
$ws.Range("D29").Value = '6.12'
$ws.Range("E29").Value = '  +1.28%  '
$ws.Range("D30").Value = '8.30'
$ws.Range("E30").Value = '  +16.39%  '
$ws.Range("D31").Value = '9.32'
$ws.Range("E31").Value = '  +3.91%  '
$ws.Range("D32").Value = '31.43'
$ws.Range("E32").Value = '  +1.00%  '
$ws.Range("D33").Value = '663.07'
$ws.Range("E33").Value = '  +6.48%  '
$ws.Range("D34").Value = '12.17'
$ws.Range("E34").Value = '  -0.05%  '
$ws.Range("D35").Value = '0.113'
$ws.Range("E35").Value = '  -0.21%  '
$ws.Range("D36").Value = '63.49'
$ws.Range("E36").Value = '  +0.14%  '
$ws.Range("D37").Value = '41.86'
$ws.Range("E37").Value = '  +0.51%  '
$ws.Range("D38").Value = '0.422'
$ws.Range("E38").Value = '  +4.59%  '
$ws.Range("E39").Value = '  -0.06%  '
$ws.Range("D40").Value = '3.288.43'
$ws.Range("E40").Value = '  +8.85%  '
$ws.Range("D41").Value = '0.0₃0760'
$ws.Range("E41").Value = '  -0.79%  '
$ws.Range("D42").Value = '3.16'
$ws.Range("E42").Value = '  +11.59%  '
$ws.Range("E43").Value = '  +2.17%  '
$ws.Range("D44").Value = '2.76'
$ws.Range("E44").Value = '  +8.97%  '
$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").Value = '0.999'
$ws.Range("E45").Value = '  -0.05%  '
$ws.Range("B46").Value = 'dogwifhat'
$ws.Range("C46").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D46").Value = '2.96'
$ws.Range("E46").Value = '  +26.58%  '
$ws.Range("D47").Value = '0.0416'
$ws.Range("E47").Value = '  +1.62%  '
$ws.Range("D48").Value = '2.75'
$ws.Range("E48").Value = '  +10.98%  '
$ws.Range("D49").Value = '8.78'
$ws.Range("E49").Value = '  +2.25%  '
$ws.Range("D50").Value = '0.131'
$ws.Range("E50").Value = '  -0.23%  '
$ws.Range("D51").Value = '3.08'
$ws.Range("E51").Value = '  -1.59%  '

# Restore default (Normal) style on the touched range so no stray
# number-format/style artifacts remain on these text cells.
$ws.Range("D2:E51").Style = "Normal"
